$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data values in B2:P6 are stored as literal text (e.g. "$108,035.92",
# "21.52%", "2,626") rather than numbers, so a plain Range.Value assignment
# would let Excel auto-parse these dollar/percent/thousand-separator-looking
# strings into real numbers. Force the whole data range to Text first so the
# new strings are kept verbatim, then clear formatting again afterwards so the
# cells end up unstyled, matching the rest of the workbook.
$dataRange = $ws.Range("B2:P6")
$dataRange.NumberFormat = "@"

$ws.Range('B2').Value = '$109,799.32'
$ws.Range('C2').Value = '$-44,673.16'
$ws.Range('D2').Value = '$-16,366.10'
$ws.Range('E2').Value = '$-20,540.63'
$ws.Range('G2').Value = '$-2,301.96'
$ws.Range('H2').Value = '$-2,301.96'
$ws.Range('J2').Value = '$-3,034.65'
$ws.Range('K2').Value = '$-10.59'
$ws.Range('L2').Value = '$-1,994.42'
$ws.Range('M2').Value = '$-997.20'
$ws.Range('N2').Value = '2,666'
$ws.Range('O2').Value = '$21,875.03'
$ws.Range('P2').Value = '19.92%'
$ws.Range('B3').Value = '$2,060.64'
$ws.Range('C3').Value = '$-441.25'
$ws.Range('D3').Value = '$-309.43'
$ws.Range('E3').Value = '$-740.32'
$ws.Range('G3').Value = '$-18.91'
$ws.Range('I3').Value = '$-0.55'
$ws.Range('J3').Value = '$-172.58'
$ws.Range('K3').Value = '$-0.18'
$ws.Range('L3').Value = '$-70.71'
$ws.Range('M3').Value = '$-17.68'
$ws.Range('N3').Value = '190'
$ws.Range('O3').Value = '$378.65'
$ws.Range('P3').Value = '18.38%'
$ws.Range('B4').Value = '$47.98'
$ws.Range('C4').Value = '$-28.01'
$ws.Range('D4').Value = '$-7.20'
$ws.Range('E4').Value = '$-8.67'
$ws.Range('J4').Value = '$-0.20'
$ws.Range('N4').Value = '2'
$ws.Range('O4').Value = '$3.90'
$ws.Range('P4').Value = '8.13%'
$ws.Range('B5').Value = '$29,458.34'
$ws.Range('C5').Value = '$-15,049.47'
$ws.Range('D5').Value = '$-4,644.32'
$ws.Range('E5').Value = '$-821.30'
$ws.Range('G5').Value = '$-336.89'
$ws.Range('I5').Value = '$-5.61'
$ws.Range('J5').Value = '$-801.23'
$ws.Range('K5').Value = '$-2.85'
$ws.Range('L5').Value = '$-686.19'
$ws.Range('M5').Value = '$-447.55'
$ws.Range('N5').Value = '132'
$ws.Range('O5').Value = '$7,686.01'
$ws.Range('P5').Value = '26.09%'
$ws.Range('J6').Value = '$-1,346.00'
$ws.Range('K6').Value = '$-26.31'
$ws.Range('L6').Value = '$-3,773.79'
$ws.Range('M6').Value = '$-3,773.79'
$ws.Range('O6').Value = '$58,111.79'
$ws.Range('P6').Value = '21.38%'

$dataRange.ClearFormats()
